$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The counterbalance sheet stores, per row, which stimulus ("stim_a" /
# "stim_d") goes with column H vs column I. For rows 3, 4, 7, 8 and 11 the
# H/I assignment needs to be flipped (columns S/T recompute automatically
# since they reference H/I via shared formulas).
$rowsToSwap = @(3, 4, 7, 8, 11)

foreach ($r in $rowsToSwap) {
    $hCell = $ws.Cells.Item($r, 8)   # column H
    $iCell = $ws.Cells.Item($r, 9)   # column I

    $hVal = $hCell.Value2
    $iVal = $iCell.Value2

    $hCell.Value2 = $iVal
    $iCell.Value2 = $hVal
}

# Update the selection to match the author's final cursor position.
$ws.Range("I15").Select() | Out-Null
